$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45913, "TV", "Spend", 126),
    @(3, 45948, "Radio", "Spend", 86),
    @(4, 45934, "Radio", "Spend", 91),
    @(5, 45906, "TV", "Spend", 187),
    @(6, 45955, "TV", "Spend", 52),
    @(7, 45927, "Radio", "Spend", 92),
    @(8, 45920, "Radio", "Spend", 85),
    @(9, 45962, "Radio", "Spend", 105),
    @(10, 45927, "TV", "Spend", 85),
    @(11, 45934, "Radio", "Spend", 91),
    @(12, 45920, "Radio", "GRPs", 5),
    @(13, 45955, "Radio", "GRPs", 7),
    @(14, 45920, "TV", "GRPs", 8),
    @(15, 45920, "TV", "GRPs", 8),
    @(16, 45913, "TV", "GRPs", 9),
    @(17, 45941, "Radio", "Spend", 92),
    @(18, 45941, "TV", "Spend", 68),
    @(19, 45927, "TV", "GRPs", 10),
    @(20, 45962, "TV", "GRPs", 2),
    @(21, 45906, "Radio", "GRPs", 5),
    @(22, 45941, "Radio", "GRPs", 5),
    @(23, 45934, "TV", "Spend", 82),
    @(24, 45962, "Radio", "GRPs", 9),
    @(25, 45920, "TV", "Spend", 198),
    @(26, 45962, "TV", "GRPs", 2),
    @(27, 45913, "Radio", "GRPs", 3),
    @(28, 45934, "Radio", "GRPs", 7),
    @(29, 45913, "Radio", "Spend", 65),
    @(30, 45955, "Radio", "Spend", 170),
    @(31, 45920, "Radio", "Spend", 85),
    @(32, 45962, "Radio", "GRPs", 9),
    @(33, 45934, "Radio", "GRPs", 7),
    @(34, 45941, "Radio", "Spend", 92),
    @(35, 45948, "TV", "GRPs", 8),
    @(36, 45913, "TV", "Spend", 126),
    @(37, 45927, "TV", "GRPs", 10),
    @(38, 45955, "TV", "Spend", 52),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
